# Updated cryptos list on Thu Oct  5 22:00:16 UTC 2023 with GitHub Actions
# Refresh price / volume(1h) figures and re-sync a few coin rows.
#
# Note: several Price values look numeric (e.g. "211.02") but must stay
# plain text (matching the original inlineStr cells), so they are entered
# with a leading apostrophe (forces text entry, like a user typing
# '211.02 in Excel) and then the cell style is reset back to Normal so no
# stray NumberFormat/quote-prefix style lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.490.36'
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('D3').Value = '1.616.21'
$ws.Range('E3').Value = '  -1.77%  '
$ws.Range('D5').Value = "'211.02"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('E6').Value = '  -1.28%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = "'22.79"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.36%  '
$ws.Range('D9').Value = "'0.261"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E10').Value = '  -0.30%  '
$ws.Range('D11').Value = "'0.0885"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.74%  '
$ws.Range('D12').Value = '1.845.75'
$ws.Range('E12').Value = '  -1.75%  '
$ws.Range('D13').Value = '1.617.25'
$ws.Range('E13').Value = '  -1.73%  '
$ws.Range('E14').Value = '  -0.29%  '
$ws.Range('E15').Value = '  -2.50%  '
$ws.Range('D16').Value = "'64.95"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.03%  '
$ws.Range('D17').Value = '27.470.64'
$ws.Range('E17').Value = '  -0.81%  '
$ws.Range('D18').Value = "'231.16"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('E19').Value = '  -1.04%  '
$ws.Range('E20').Value = '  -2.05%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('E22').Value = '  -0.95%  '
$ws.Range('E23').Value = '  +0.77%  '
$ws.Range('D24').Value = "'2.08"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.65%  '
$ws.Range('D25').Value = "'150.79"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.63%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = "'6.84"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.93%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = "'0.111"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.05%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('E30').Value = '  -1.05%  '
$ws.Range('E31').Value = '  -0.81%  '
$ws.Range('E32').Value = '  -1.11%  '
$ws.Range('D33').Value = '1.464.84'
$ws.Range('E33').Value = '  +1.34%  '
$ws.Range('E34').Value = '  -2.99%  '
$ws.Range('E35').Value = '  -4.01%  '
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('D37').Value = "'0.941"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.54%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = "'0.559"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.45%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = "'0.0167"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.52%  '
$ws.Range('E40').Value = '  -2.95%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').Value = "'67.96"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.00%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = "'0.988"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.38%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').Value = "'2.20"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.28%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = "'5.26"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.38%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.756.37'
$ws.Range('E46').Value = '  -1.78%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = "'1.71"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.25%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = "'86.61"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0105'
$ws.Range('E49').Value = '  -2.31%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = "'0.101"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.55%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = "'7.67"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.87%  '